# Edit script generated to match target diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert 10 new rows before row 99 to make room for additional data rows
# (old row 99 becomes row 109, old rows 104/105 become 114/115)
$ws.Rows("99:108").Insert(-4121)

# Step 2: Copy formatting from row 98 (a standard data row) into the newly inserted rows
# Restrict to columns B:J so we don't bloat the used range across the whole row
$ws.Range("B98:J98").Copy()
$ws.Range("B99:J108").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Step 3: Write full data table for rows 16-109 (B,C,D,E,F,G columns)
$data = @(
    @(16,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2507",31249,781242),
    @(17,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2506",31249,781242),
    @(18,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2505",31249,781242),
    @(19,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2504",31249,781242),
    @(20,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2503",31249,781242),
    @(21,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2502",31249,781242),
    @(22,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2501",31249,781242),
    @(23,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2412",31249,781242),
    @(24,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2411",31249,781242),
    @(25,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2410",31249,781242),
    @(26,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2409",31249,781242),
    @(27,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2408",31249,781242),
    @(28,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2407",31249,781242),
    @(29,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2406",31249,781242),
    @(30,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2405",31249,781242),
    @(31,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2404",31249,781242),
    @(32,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2403",31249,781242),
    @(33,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2402",31249,781242),
    @(34,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2401",31249,781242),
    @(35,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2312",31249,781242),
    @(36,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2311",31249,781242),
    @(37,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2310",31249,781242),
    @(38,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2309",31249,781242),
    @(39,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2308",31249,781242),
    @(40,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2307",31249,781242),
    @(41,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2306",31249,781242),
    @(42,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2305",31249,781242),
    @(43,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2304",31249,781242),
    @(44,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2303",31249,781242),
    @(45,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2302",31249,781242),
    @(46,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2301",31249,781242),
    @(47,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2212",31249,781242),
    @(48,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2211",31249,781242),
    @(49,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2210",31249,781242),
    @(50,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2209",31249,781242),
    @(51,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2208",31249,781242),
    @(52,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2207",31249,781242),
    @(53,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2206",31249,781242),
    @(54,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2205",31249,781242),
    @(55,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2204",31249,781242),
    @(56,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2203",31249,781242),
    @(57,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2202",31249,781242),
    @(58,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2201",31249,781242),
    @(59,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2112",31249,781242),
    @(60,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2111",31249,781242),
    @(61,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2110",31249,781242),
    @(62,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2109",31249,781242),
    @(63,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2108",31249,781242),
    @(64,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2107",31249,781242),
    @(65,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2106",31249,781242),
    @(66,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2105",31249,781242),
    @(67,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2104",31249,781242),
    @(68,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2103",31249,781242),
    @(69,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2102",31249,781242),
    @(70,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2101",31249,781242),
    @(71,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2012",31249,781242),
    @(72,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2011",31249,781242),
    @(73,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2010",31249,781242),
    @(74,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2009",31249,781242),
    @(75,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2008",31249,781242),
    @(76,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2007",31249,781242),
    @(77,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2006",31249,781242),
    @(78,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2005",31249,781242),
    @(79,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2004",31249,781242),
    @(80,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2003",31249,781242),
    @(81,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2002",31249,781242),
    @(82,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","2001",31249,781242),
    @(83,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","1912",31249,781242),
    @(84,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","1911",31249,781242),
    @(85,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","1910",31249,781242),
    @(86,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","1909",31249,781242),
    @(87,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","1908",31249,781242),
    @(88,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","1907",31249,781242),
    @(89,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","1906",31249,781242),
    @(90,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","1905",31249,781242),
    @(91,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","1904",31249,781242),
    @(92,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","1903",31249,781242),
    @(93,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","1902",31249,781242),
    @(94,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","1901",31249,781242),
    @(95,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","1812",31249,781242),
    @(96,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","1811",31249,781242),
    @(97,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","1810",31249,781242),
    @(98,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","1809",31249,781242),
    @(99,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","1808",31249,781242),
    @(100,"CC","1115856975","AMPARO ASTRID GARCIA ESCOBAR","1807",31249,781242),
    @(101,"CC","45550109","BERLYS MARIA ARRAZOLA OSORIO","2507",88711,2217767),
    @(102,"CC","1047424156","KIMBERLY JOHANA WATTS DIAZ","2507",56940,1423500),
    @(103,"CC","1047398676","EDITH CAROLA GARCIA CASTELLAR","2507",113856,2846391),
    @(104,"CC","1047405655","JONATHAN RAFAEL MAZO LOPEZ","2507",52000,1300000),
    @(105,"CC","1126787163","ESTEFANY ALEJANDRA MAZO CAÑAS","2507",56940,1423500),
    @(106,"CC","1126787163","ESTEFANY ALEJANDRA MAZO CAÑAS","1810",56940,1423500),
    @(107,"CC","92446752","JUAN BALLESTERO MURILLO","2507",56940,1423500),
    @(108,"CC","1103219791","MANUEL DAVID PEREZ ORTIZ","2507",56940,1423500),
    @(109,"CC","1047395467","MICHAEL JOSE JIMENEZ LOPEZ","2507",56940,1423500)
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells($r, 2).Value2 = $item[1]   # B: Tipo Doc
    $ws.Cells($r, 3).Value2 = $item[2]   # C: N Doc
    $ws.Cells($r, 4).Value2 = $item[3]   # D: Nombre
    $ws.Cells($r, 5).Value2 = $item[4]   # E: Periodo Mora
    $ws.Cells($r, 6).Value2 = $item[5]   # F: Valor Mora
    $ws.Cells($r, 7).Value2 = $item[6]   # G: Salario Basico
}

# Step 4: Update summary cells
$ws.Range("E11").Value2 = 3252372   # Valor Mora total
$ws.Range("C13").Value2 = 9         # Cant. Trabajadores
$ws.Range("F13").Value2 = 85        # Cant. Periodos
